# Apply the target changes to the workbook:
#  1. Add a new blank worksheet named "Sheet3" at the end.
#  2. On the "test_suite" sheet, change A3 from "loginTest" to "ProductPage"
#     (with an underlined-font style), and add two new rows:
#       A4 = "loginTest"   / B4 = "Y"
#       A5 = "productPage" (underlined-font style) / B5 = "Y"

$wb = $excel.ActiveWorkbook

$wsSuite = $wb.Worksheets.Item("test_suite")

# --- Update existing row 3: loginTest -> ProductPage -------------------
$wsSuite.Range("A3").Value = "ProductPage"
$wsSuite.Range("A3").Font.Underline = 2
$wsSuite.Range("A3").Font.Name = "Arial"
$wsSuite.Range("A3").Font.Size = 10
$wsSuite.Range("A3").Font.ColorIndex = 1

$wsSuite.Range("B3").Value = "Y"

# --- Add row 4: loginTest / Y ------------------------------------------
$wsSuite.Range("A4").Value = "loginTest"
$wsSuite.Range("B4").Value = "Y"

# --- Add row 5: productPage / Y -----------------------------------------
$wsSuite.Range("A5").Value = "productPage"
$wsSuite.Range("A5").Font.Underline = 2
$wsSuite.Range("A5").Font.Name = "Arial"
$wsSuite.Range("A5").Font.Size = 10
$wsSuite.Range("A5").Font.ColorIndex = 1

$wsSuite.Range("B5").Value = "Y"

# --- Add the new blank "Sheet3" worksheet at the end ---------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "Sheet3"

# Keep "test_suite" as the active/selected sheet (matches original state).
$wsSuite.Activate()
$wsSuite.Range("C14").Select()
